# feat: add 2022-Q3 data
#
# 1) Insert a new "2022-Q3" worksheet right after "总计" (i.e. as the new
#    2nd sheet, pushing every later quarter sheet back by one position).
# 2) Fill it with that quarter's per-fund holding detail (same shape as the
#    other quarter sheets: header row + 7 fund rows).
# 3) Prepend a matching summary row to the "总计" sheet, shifting the
#    existing history rows down by one and renumbering the index column.

$wb = $excel.ActiveWorkbook
$summary = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# Step 1: new "2022-Q3" sheet, inserted right after "总计"
# ---------------------------------------------------------------------------
$q3 = $wb.Worksheets.Add($null, $summary)
$q3.Name = "2022-Q3"

function Set-TextCell($ws, $addr, $text) {
    # Force "text" storage even for numeric-looking strings (fund codes like
    # "012640", ratio strings like "2.61") so leading zeros / exact strings
    # survive, matching the source sheets where columns B..G are text.
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $text
}

# Header row (values first; style gets fixed up below)
Set-TextCell $q3 "B1" "基金代码"
Set-TextCell $q3 "C1" "基金名称"
Set-TextCell $q3 "D1" "基金规模"
Set-TextCell $q3 "E1" "股票总仓位"
Set-TextCell $q3 "F1" "仓位占比"
Set-TextCell $q3 "G1" "持有市值(亿元)"
Set-TextCell $q3 "H1" "仓位排名"

$q3Rows = @(
    @("012640", "鹏华稳健鸿利一年持有期混合A", "2.61", "92.98", "2.82", "0.0736", 10),
    @("008134", "鹏华优选价值股票", "1.80", "92.72", "2.67", "0.0481", 10),
    @("011574", "鹏华领航一年持有期混合A", "1.20", "92.84", "3.25", "0.0390", 9),
    @("012010", "富国泰享回报6个月持有期混合A", "6.39", "29.75", "0.61", "0.0390", 10),
    @("011575", "鹏华领航一年持有期混合C", "0.91", "92.84", "3.25", "0.0296", 9),
    @("012641", "鹏华稳健鸿利一年持有期混合C", "0.10", "92.98", "2.82", "0.0028", 10),
    @("012011", "富国泰享回报6个月持有期混合C", "0.05", "29.75", "0.61", "0.0003", 10)
)

for ($i = 0; $i -lt $q3Rows.Count; $i++) {
    $r = $i + 2
    $data = $q3Rows[$i]
    $q3.Range("A$r").Value = $i
    Set-TextCell $q3 "B$r" $data[0]
    Set-TextCell $q3 "C$r" $data[1]
    Set-TextCell $q3 "D$r" $data[2]
    Set-TextCell $q3 "E$r" $data[3]
    Set-TextCell $q3 "F$r" $data[4]
    Set-TextCell $q3 "G$r" $data[5]
    $q3.Range("H$r").Value = $data[6]
}

# The NumberFormat="@" trick above is only needed transiently (so Excel
# stores "012640" / "2.61" as text instead of re-parsing them as numbers);
# it leaves every touched cell on a throwaway "text" style. Restore the
# real look of each area by pasting formats from cells on "总计" that
# already carry the right (untouched) style:
#   - B1:H1 (header)      <- bold+border centered header style
#   - A2:A8 (index column) <- same centered style
#   - B2:G8 (plain data)   <- no explicit style at all
$summary.Range("B1").Copy()
$q3.Range("B1:H1").PasteSpecial(-4122)   # xlPasteFormats
$summary.Range("A2").Copy()
$q3.Range("A2:A8").PasteSpecial(-4122)   # xlPasteFormats
$summary.Range("C2").Copy()
$q3.Range("B2:G8").PasteSpecial(-4122)   # xlPasteFormats

# ---------------------------------------------------------------------------
# Step 2: push the "总计" history down one row and insert the 2022-Q3 totals
# ---------------------------------------------------------------------------
$lastRow = 8
for ($r = $lastRow; $r -ge 2; $r--) {
    $src = $summary.Range("A" + $r + ":D" + $r)
    $dst = $summary.Range("A" + ($r + 1) + ":D" + ($r + 1))
    $src.Copy($dst)
}

for ($r = 2; $r -le ($lastRow + 1); $r++) {
    $summary.Range("A$r").Value = $r - 2
}

$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 7
$summary.Range("D2").Value = 0.23
